$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1850
$ws.Range("I70").Value = 700
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 2100
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -1830
$ws.Range("N70").Value = -9540

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1850
$ws.Range("I73").Value = 700
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 2100
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -1164
$ws.Range("N73").Value = -10872

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 11454.77
$ws.Range("I132").Value = 13526.944
$ws.Range("K132").Value = 40580.83199999999
$ws.Range("M132").Value = -38050.83199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2237.5488
$ws.Range("I138").Value = 1733.3478
$ws.Range("K138").Value = 5200.0434
$ws.Range("M138").Value = -60.04340000000047

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 3167.389
$ws.Range("I5").Value = 448.33334
$ws.Range("J5").Value = 5886.4443
$ws.Range("K5").Value = 448.33334
$ws.Range("L5").Value = 5886.4443
$ws.Range("M5").Value = -336.33334
$ws.Range("N5").Value = -6110.4443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 85845
$ws.Range("J7").Value = 85845
$ws.Range("L7").Value = 85845
$ws.Range("N7").Value = -86073

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 139572.47
$ws.Range("I32").Value = 148758.89
$ws.Range("K32").Value = 148758.89
$ws.Range("M32").Value = -148471.89

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 29412986
$ws.Range("I97").Value = 1030.8966
$ws.Range("K97").Value = 1030.8966
$ws.Range("M97").Value = -534.8966

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1899.2
$ws.Range("J122").Value = 2599.6667
$ws.Range("L122").Value = 7799.000100000001
$ws.Range("N122").Value = -12699.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 79999.336
$ws.Range("J138").Value = 79999.336
$ws.Range("L138").Value = 79999.336
$ws.Range("N138").Value = -90279.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 3167.389
$ws.Range("I4").Value = 448.33334
$ws.Range("J4").Value = 5886.4443
$ws.Range("K4").Value = 448.33334
$ws.Range("L4").Value = 5886.4443
$ws.Range("M4").Value = -333.33334
$ws.Range("N4").Value = -6116.4443

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 105000
$ws.Range("J59").Value = 105000
$ws.Range("L59").Value = 105000
$ws.Range("N59").Value = -106694

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 31420.385
$ws.Range("J82").Value = 40855.445
$ws.Range("L82").Value = 40855.445
$ws.Range("N82").Value = -41621.445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 31420.385
$ws.Range("J85").Value = 40855.445
$ws.Range("L85").Value = 40855.445
$ws.Range("N85").Value = -43507.445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2126.3242
$ws.Range("I94").Value = 1933.5714
$ws.Range("K94").Value = 1933.5714
$ws.Range("M94").Value = -1482.5714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 42847.6
$ws.Range("I135").Value = 16709
$ws.Range("J135").Value = 49382.25
$ws.Range("K135").Value = 16709
$ws.Range("L135").Value = 49382.25
$ws.Range("M135").Value = -11639
$ws.Range("N135").Value = -59522.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 866.2857
$ws.Range("I16").Value = 862.5
$ws.Range("J16").Value = 871.3333
$ws.Range("K16").Value = 862.5
$ws.Range("L16").Value = 871.3333
$ws.Range("M16").Value = -575.5
$ws.Range("N16").Value = -1445.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 6882.1577
$ws.Range("I94").Value = 14160.5
$ws.Range("J94").Value = 1588.8182
$ws.Range("K94").Value = 14160.5
$ws.Range("L94").Value = 1588.8182
$ws.Range("M94").Value = -13709.5
$ws.Range("N94").Value = -2490.8182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3580.7778
$ws.Range("I99").Value = 4035.5
$ws.Range("J99").Value = 2671.3333
$ws.Range("K99").Value = 4035.5
$ws.Range("L99").Value = 2671.3333
$ws.Range("M99").Value = -2537.5
$ws.Range("N99").Value = -5667.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 866.2857
$ws.Range("I113").Value = 862.5
$ws.Range("J113").Value = 871.3333
$ws.Range("K113").Value = 862.5
$ws.Range("L113").Value = 871.3333
$ws.Range("M113").Value = 1307.5
$ws.Range("N113").Value = -5211.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 8017.8965
$ws.Range("I122").Value = 1842.0454
$ws.Range("K122").Value = 5526.1362
$ws.Range("M122").Value = -3076.1362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3580.7778
$ws.Range("I126").Value = 4035.5
$ws.Range("J126").Value = 2671.3333
$ws.Range("K126").Value = 12106.5
$ws.Range("L126").Value = 8013.999899999999
$ws.Range("M126").Value = -9636.5
$ws.Range("N126").Value = -12953.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2247.5881
$ws.Range("I134").Value = 1923.1786
$ws.Range("K134").Value = 5769.5358
$ws.Range("M134").Value = -3234.5358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8616.412
$ws.Range("I3").Value = 3184.875
$ws.Range("J3").Value = 13444.444
$ws.Range("K3").Value = 9554.625
$ws.Range("L3").Value = 40333.33199999999
$ws.Range("M3").Value = -9442.625
$ws.Range("N3").Value = -40557.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 13300
$ws.Range("I94").Value = 6500
$ws.Range("K94").Value = 19500
$ws.Range("M94").Value = -18824

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 717223.3
$ws.Range("I129").Value = 1112927.8
$ws.Range("J129").Value = 4955.2
$ws.Range("K129").Value = 3338783.4
$ws.Range("L129").Value = 14865.6
$ws.Range("M129").Value = -3333783.4
$ws.Range("N129").Value = -24865.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 80899.16
$ws.Range("I122").Value = 128461.25
$ws.Range("J122").Value = 4799.8
$ws.Range("K122").Value = 385383.75
$ws.Range("L122").Value = 14399.4
$ws.Range("M122").Value = -382933.75
$ws.Range("N122").Value = -19299.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2155.9688
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2155.9688
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2155.9688
$ws.Range("N22").Value = -2745.9688
$ws.Range("M22").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2155.9688
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2155.9688
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2155.9688
$ws.Range("N27").Value = -2369.9688
$ws.Range("M27").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5290.8
$ws.Range("I46").Value = 954.5
$ws.Range("J46").Value = 6374.875
$ws.Range("K46").Value = 954.5
$ws.Range("L46").Value = 6374.875
$ws.Range("M46").Value = -766.5
$ws.Range("N46").Value = -6750.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1498.7858
$ws.Range("I55").Value = 2074.7778
$ws.Range("K55").Value = 2074.7778
$ws.Range("M55").Value = -1901.7778

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H134").Value = 90000
$ws.Range("J134").Value = 90000
$ws.Range("L134").Value = 90000
$ws.Range("N134").Value = -100140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 1000
$ws.Range("J26").Value = 1000
$ws.Range("L26").Value = 1000
$ws.Range("N26").Value = -1586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 70997.5
$ws.Range("I122").Value = 466.16666
$ws.Range("J122").Value = 282591.5
$ws.Range("K122").Value = 1398.49998
$ws.Range("L122").Value = 847774.5
$ws.Range("M122").Value = 1051.50002
$ws.Range("N122").Value = -852674.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2575.2666
$ws.Range("I126").Value = 2657.5454
$ws.Range("J126").Value = 2349
$ws.Range("K126").Value = 7972.6362
$ws.Range("L126").Value = 7047
$ws.Range("M126").Value = -5502.6362
$ws.Range("N126").Value = -11987

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 65000
$ws.Range("J135").Value = 65000
$ws.Range("L135").Value = 65000
$ws.Range("N135").Value = -75140

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4384.533
$ws.Range("I136").Value = 4619.875
$ws.Range("J136").Value = 4115.5713
$ws.Range("K136").Value = 13859.625
$ws.Range("L136").Value = 12346.7139
$ws.Range("M136").Value = -11309.625
$ws.Range("N136").Value = -17446.7139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H137").Value = 60000
$ws.Range("J137").Value = 60000
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200
